$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header counts (number of observations per condition) changed
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 (Lichtwark) values replaced
$ws.Range("B2").Value = 260.12934017581102
$ws.Range("C2").Value = 302.87295104901125
$ws.Range("D2").Value = 255.41594220805601
$ws.Range("E2").Value = 306.94565078381743

# Row 3 values replaced
$ws.Range("B3").Value = 249.29294584859031
$ws.Range("C3").Value = 302.93054679703334
$ws.Range("D3").Value = 250.41660847738984
$ws.Range("E3").Value = 307.96892918008837

# Selection narrowed to the edited block
$ws.Range("B1:E3").Select()
